# Apply data update to the "Inscricoes" sheet (Table1) as described in the
# commit "Data update using git".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 5: Inscritos 5 -> 6
$ws.Range("E5").Value = 6

# Row 19: Inscritos 59 -> 60
$ws.Range("E19").Value = 60

# Row 28: Inscritos 19 -> 20, Pagos 15 -> 16, Inscricoes homologadas 17 -> 18
$ws.Range("E28").Value = 20
$ws.Range("F28").Value = 16
$ws.Range("H28").Value = 18

# Row 32: Inscritos 21 -> 20
$ws.Range("E32").Value = 20

# Row 34: Inscritos 23 -> 24
$ws.Range("E34").Value = 24

# Row 36: Inscritos 106 -> 107, Pagos 46 -> 47, Inscricoes homologadas 78 -> 79
$ws.Range("E36").Value = 107
$ws.Range("F36").Value = 47
$ws.Range("H36").Value = 79

# Row 51: Inscritos 12 -> 13
$ws.Range("E51").Value = 13

# Row 79: Inscritos 38 -> 39
$ws.Range("E79").Value = 39

# Row 89: Inscritos 43 -> 45
$ws.Range("E89").Value = 45
